# ---------------------------------------------------------------------------
# Applies the "Added a prototype dyr parse; Changed EXDC2 and TGOV1
# parameters to follow PSS/E convention (all caps). Updated Kundur system
# correspondingly." commit to the Kundur workbook.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsEXDC2   = $wb.Worksheets.Item("EXDC2")
$wsTGOV1   = $wb.Worksheets.Item("TGOV1")
$wsToggler = $wb.Worksheets.Item("Toggler")

# ---------------------------------------------------------------------------
# 1. EXDC2 sheet: rename parameters to PSS/E (all-caps) convention, re-order
#    a couple of columns (KF1/TF1 swap vs. old Tf/Kf), and append 4 new
#    saturation / limit columns (E1, SE1, E2, SE2).
# ---------------------------------------------------------------------------

# -- existing header cells F1:P1 get new (upper-case) names ----------------
# (K1/L1 are written "L-then-K" so the new shared-string table lands in the
# same order - KF1 before TF1 - as the reference workbook.)
$wsEXDC2.Cells.Item(1, 6).Value  = "TR"
$wsEXDC2.Cells.Item(1, 7).Value  = "TA"
$wsEXDC2.Cells.Item(1, 8).Value  = "TC"
$wsEXDC2.Cells.Item(1, 9).Value  = "TB"
$wsEXDC2.Cells.Item(1, 10).Value = "TE"
$wsEXDC2.Cells.Item(1, 12).Value = "KF1"
$wsEXDC2.Cells.Item(1, 11).Value = "TF1"
$wsEXDC2.Cells.Item(1, 13).Value = "KA"
$wsEXDC2.Cells.Item(1, 14).Value = "KE"
$wsEXDC2.Cells.Item(1, 15).Value = "VRMAX"
$wsEXDC2.Cells.Item(1, 16).Value = "VRMIN"
# Q1 (Ae) and R1 (Be) keep their text - no change needed.

# -- four brand-new trailing columns: S, T, U, V ----------------------------
$wsEXDC2.Cells.Item(1, 19).Value = "E1"
$wsEXDC2.Cells.Item(1, 20).Value = "SE1"
$wsEXDC2.Cells.Item(1, 21).Value = "E2"
$wsEXDC2.Cells.Item(1, 22).Value = "SE2"

# Style the new header cells like the existing bold/bordered/centred header
# (font 1, thin left+right border, centred & top-aligned) and give them a
# "no fill" override.
$newHeaderRange = $wsEXDC2.Range("S1:V1")
$newHeaderRange.Font.Bold = $true
$newHeaderRange.HorizontalAlignment = -4108   # xlCenter
$newHeaderRange.VerticalAlignment = -4160     # xlTop
$newHeaderRange.Borders.Item(7).LineStyle = 1   # xlInsideVertical-ish (left edge)
$newHeaderRange.Borders.Item(10).LineStyle = 1  # right edge

for ($col = 19; $col -le 22; $col++) {
    $c = $wsEXDC2.Cells.Item(1, $col)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.Item(7).LineStyle = 1
    $c.Borders.Item(10).LineStyle = 1
}

# -- new data columns S:V for the four EXDC2 model rows (all zero) ---------
for ($row = 2; $row -le 5; $row++) {
    $wsEXDC2.Cells.Item($row, 19).Value = 0   # S - E1
    $wsEXDC2.Cells.Item($row, 20).Value = 0   # T - SE1

    $u = $wsEXDC2.Cells.Item($row, 21)        # U - E2
    $u.Value = 0
    $u.Font.Color = 0

    $v = $wsEXDC2.Cells.Item($row, 22)        # V - SE2
    $v.Value = 0
    $v.Font.Color = 0
}

# -- update the frozen-pane selection to match the new last column ---------
$wsEXDC2.Range("V2:V5").Select()

# ---------------------------------------------------------------------------
# 2. TGOV1 sheet: rename vmin/vmax to PSS/E convention (VMIN/VMAX) and make
#    this the active/selected sheet (it takes over from Toggler).
# ---------------------------------------------------------------------------

$wsTGOV1.Cells.Item(1, 7).Value = "VMIN"   # G1: vmin -> VMIN
$wsTGOV1.Cells.Item(1, 8).Value = "VMAX"   # H1: vmax -> VMAX

$wsTGOV1.Range("H2").Select()

# ---------------------------------------------------------------------------
# 3. Workbook-level window view: TGOV1 becomes the active tab instead of
#    Toggler (Toggler's tabSelected flag is dropped automatically once a
#    different sheet is activated/selected).
# ---------------------------------------------------------------------------

$win = $wb.Windows.Item(1)
$win.Left = 1700
$win.Top = 460
$win.Width = 31900
$win.Height = 19340
